# Edit summary:
#   1. Slide 6's table ("Sources of finance") gets a different built-in
#      PowerPoint table style applied (GUID changes from
#      {3140676B-44D8-4EDB-852F-D1187D1AF1FC} to
#      {6A70DBDE-FF50-44F7-9C50-4F70CF5F644D}).
#   2. The deck's theme palette is swapped from the custom "Integral"
#      colour scheme back to the stock Office colour scheme (the 12
#      theme colours that PowerPoint ships with by default).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Re-style the table on slide 6 via the real Table object model —
#    Table.Style can't be assigned directly, PowerPoint requires
#    Table.ApplyStyle("{GUID}").
# ---------------------------------------------------------------------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle("{6A70DBDE-FF50-44F7-9C50-4F70CF5F644D}")
        }
    }
}

# ---------------------------------------------------------------------
# 2) Restore the stock "Office" theme colours on the deck's theme.
#    PowerShell COM-interop has no RGB() helper (that's VBA-only), so
#    pack R/G/B into the OLE COLORREF (0x00BBGGRR) by hand.
# ---------------------------------------------------------------------
function ColorRef([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$officeColors = @(
    (ColorRef 0x00 0x00 0x00),   # 1  dk1
    (ColorRef 0xFF 0xFF 0xFF),   # 2  lt1
    (ColorRef 0x44 0x54 0x6A),   # 3  dk2
    (ColorRef 0xE7 0xE6 0xE6),   # 4  lt2
    (ColorRef 0x5B 0x9B 0xD5),   # 5  accent1
    (ColorRef 0xED 0x7D 0x31),   # 6  accent2
    (ColorRef 0xA5 0xA5 0xA5),   # 7  accent3
    (ColorRef 0xFF 0xC0 0x00),   # 8  accent4
    (ColorRef 0x44 0x72 0xC4),   # 9  accent5
    (ColorRef 0x70 0xAD 0x47),   # 10 accent6
    (ColorRef 0x05 0x63 0xC1),   # 11 hlink
    (ColorRef 0x95 0x4F 0x72)    # 12 folHlink
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($k = 1; $k -le $officeColors.Count; $k++) {
    $themeColors.Item($k).RGB = $officeColors[$k - 1]
}
